$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 3787
    3  = 3963
    4  = 3993
    5  = 4225
    6  = 4225
    7  = 4270
    8  = 4436
    9  = 4436
    10 = 4436
    11 = 4436
    12 = 4461
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 3).Value = $values[$row]
}
